$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-11: A (Sending cluster), D (Target cluster), and numeric columns E..T.
# B is always "Fn1" and C is always "Tnfrsf11b" for every row (unchanged from before).
$rows = @(
    @{ A="ECs";  D="ECs";  E=3; F=1; G=61.84465033333333;  H=185.533951;        I=0.03153184209101587; J=0.03153184209101587; K=1; L=0.3333333333333333; M=0.08241233333333334; N=0.247237;   O=0.04727005612861496; P=0.04727005612861496; Q=5.096761938154112;    R=45.870857443387;    S=0.001490511945480944; T=0.001490511945480944 },
    @{ A="ECs";  D="FAPs"; E=3; F=1; G=61.84465033333333;  H=185.533951;        I=0.03153184209101587; J=0.03153184209101587; K=3; L=1;                   M=1.661024;             N=4.983072;    O=0.9527299438713851;  P=0.952729943871385;   Q=102.7254484752747;    R=924.529036277472;   S=0.03004133014553493;  T=0.03004133014553492  },
    @{ A="FAPs"; D="ECs";  E=3; F=1; G=1361.379069;        H=4084.137207;       I=0.6941067594101231;  J=0.6941067594101232;  K=1; L=0.3333333333333333; M=0.08241233333333334; N=0.247237;   O=0.04727005612861496; P=0.04727005612861496; Q=112.194425627451;     R=1009.749830647059;  S=0.03281046547656756;  T=0.03281046547656757  },
    @{ A="FAPs"; D="FAPs"; E=3; F=1; G=1361.379069;        H=4084.137207;       I=0.6941067594101231;  J=0.6941067594101232;  K=3; L=1;                   M=1.661024;             N=4.983072;    O=0.9527299438713851;  P=0.952729943871385;   Q=2261.283306706656;    R=20351.5497603599;   S=0.6612962939335556;   T=0.6612962939335556   },
    @{ A="M1";   D="ECs";  E=3; F=1; G=251.007014;         H=753.021042;        I=0.127977334927537;   J=0.1279773349275369;  K=1; L=0.3333333333333333; M=0.08241233333333334; N=0.247237;   O=0.04727005612861496; P=0.04727005612861496; Q=20.68607370677267;    R=186.174663360954;   S=0.006049495805215228; T=0.006049495805215227 },
    @{ A="M1";   D="FAPs"; E=3; F=1; G=251.007014;         H=753.021042;        I=0.127977334927537;   J=0.1279773349275369;  K=3; L=1;                   M=1.661024;             N=4.983072;    O=0.9527299438713851;  P=0.952729943871385;   Q=416.928674422336;     R=3752.358069801024;  S=0.1219278391223217;   T=0.1219278391223217   },
    @{ A="M2";   D="ECs";  E=3; F=1; G=260.0315303333334;  H=780.094591;        I=0.1325785352324417;  J=0.1325785352324417;  K=1; L=0.3333333333333333; M=0.08241233333333334; N=0.247237;   O=0.04727005612861496; P=0.04727005612861496; Q=21.42980515500745;    R=192.868246395067;   S=0.006266994801887076; T=0.006266994801887076 },
    @{ A="M2";   D="FAPs"; E=3; F=1; G=260.0315303333334;  H=780.094591;        I=0.1325785352324417;  J=0.1325785352324417;  K=3; L=1;                   M=1.661024;             N=4.983072;    O=0.9527299438713851;  P=0.952729943871385;   Q=431.9186126403947;    R=3887.267513763552;  S=0.1263115404305546;   T=0.1263115404305546   },
    @{ A="sCs";  D="ECs";  E=3; F=1; G=27.07732933333333;  H=81.231988;         I=0.01380552833888228; J=0.01380552833888228; K=1; L=0.3333333333333333; M=0.08241233333333334; N=0.247237;   O=0.04727005612861496; P=0.04727005612861496; Q=2.231505890795111;    R=20.083553017156;    S=0.0006525880994641499; T=0.0006525880994641499 },
    @{ A="sCs";  D="FAPs"; E=3; F=1; G=27.07732933333333;  H=81.231988;         I=0.01380552833888228; J=0.01380552833888228; K=3; L=1;                   M=1.661024;             N=4.983072;    O=0.9527299438713851;  P=0.952729943871385;   Q=44.97609387857067;    R=404.784844907136;   S=0.01315294023941813;  T=0.01315294023941813  }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = "Fn1"
    $ws.Cells.Item($r, 3).Value = "Tnfrsf11b"
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $r = $r + 1
}
